$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '30.204.86'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '1.861.00'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'236.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").Value = "'0.4676"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.51%  '
$ws.Range("D8").Value = "'0.2851"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.39%  '
$ws.Range("D9").Value = "'0.06524"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.35%  '
$ws.Range("E10").Value = '  +9.96%  '
$ws.Range("D11").Value = "'0.07896"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.20%  '
$ws.Range("D12").Value = "'97.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("D13").Value = '1.864.64'
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").Value = "'5.163"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.12%  '
$ws.Range("D15").Value = "'0.6794"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.31%  '
$ws.Range("D16").Value = "'279.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.16%  '
$ws.Range("D17").Value = '30.199.06'
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("E18").Value = '  +7.04%  '
$ws.Range("D19").Value = "'1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("E20").Value = '  -1.60%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = "'0.000007310"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.13%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.108.47'
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = "'6.159"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("D25").Value = "'167.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("D26").Value = "'9.227"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.06%  '
$ws.Range("D27").Value = "'19.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.88%  '
$ws.Range("D28").Value = "'1.930"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.09%  '
$ws.Range("E29").Value = '  +3.23%  '
$ws.Range("D30").Value = "'0.09732"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.88%  '
$ws.Range("D31").Value = "'4.366"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.20%  '
$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("D33").Value = "'4.046"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.35%  '
$ws.Range("D34").Value = "'0.04726"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.89%  '
$ws.Range("D35").Value = "'1.134"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.23%  '
$ws.Range("D36").Value = "'0.7075"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.26%  '
$ws.Range("D37").Value = "'2.708"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("E38").Value = '  +1.45%  '
$ws.Range("D39").Value = "'2.630"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.84%  '
$ws.Range("D40").Value = "'6.308"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.36%  '
$ws.Range("D41").Value = "'74.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.60%  '
$ws.Range("D42").Value = "'1.950"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.71%  '
$ws.Range("D43").Value = "'0.8482"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.61%  '
$ws.Range("D44").Value = "'0.4173"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.80%  '
$ws.Range("D45").Value = "'0.9997"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("D46").Value = "'103.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("D47").Value = "'969.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("D48").Value = "'7.181"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.14%  '
$ws.Range("D49").Value = "'9.274"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.90%  '
$ws.Range("D50").Value = "'34.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("D51").Value = "'0.05637"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.22%  '
